$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12/13 swap (TRON <-> WrappedEther)
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"

# Row 18/19 swap (Avalanche <-> Dai)
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

# Price (D) and Volume(1h) (E) text values - force text format to preserve exact string
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.912.89"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "1.857.30"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "305.50"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.5053"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("D8").Value = "0.3714"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "0.07122"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "0.8836"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "20.49"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "1.876.96"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").Value = "0.07545"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").Value = "5.276"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "88.59"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "0.000008345"
$ws.Range("E17").Value = "  -4.46%  "
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "14.04"
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").Value = "26.938.06"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").Value = "5.032"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").Value = "2.095.63"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("D23").Value = "10.46"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").Value = "6.448"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").Value = "1.847"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "146.93"
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("D27").Value = "17.94"
$ws.Range("E27").Value = "  -2.05%  "
$ws.Range("D28").Value = "2.084"
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("D29").Value = "112.39"
$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("D30").Value = "4.645"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").Value = "4.633"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("D32").Value = "0.09035"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "0.05098"
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("D34").Value = "3.046"
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").Value = "1.144"
$ws.Range("E35").Value = "  -7.30%  "
$ws.Range("D36").Value = "0.7232"
$ws.Range("E36").Value = "  -7.35%  "
$ws.Range("D37").Value = "0.02030"
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").Value = "3.033"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").Value = "2.458"
$ws.Range("E39").Value = "  -6.35%  "
$ws.Range("D40").Value = "1.069"
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").Value = "0.5271"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("D42").Value = "6.540"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("D43").Value = "114.91"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "8.217"
$ws.Range("E44").Value = "  -3.13%  "
$ws.Range("D45").Value = "0.1465"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "0.4580"
$ws.Range("E47").Value = "  -4.25%  "
$ws.Range("D48").Value = "9.955"
$ws.Range("E48").Value = "  -4.52%  "
$ws.Range("D49").Value = "1.549"
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("D50").Value = "36.37"
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("D51").Value = "63.68"
$ws.Range("E51").Value = "  -4.48%  "
